$wb = $excel.ActiveWorkbook
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $wsLast)
$ws.Name = "Sheet1"
$ws.Columns("A:B").ColumnWidth = 9.140625
$ws.Columns("C").ColumnWidth = 11.85546875
$ws.Columns("D").ColumnWidth = 12
